$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial for every data row (2-103).
# All of them were bumped by one day (45181 -> 45182).
$range = $ws.Range("C2:C103")
$range.Value = 45182
